# Update Efna4-Epha5.xlsx sheet with new TPM-based NATMI output (3 -> 6 data rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Efna4", "Epha5", "ECs", 2, 0.6666666666666666, 0.2988413333333333, 0.896524, 0.3632971504731247, 0.3632971504731246, 1, 0.3333333333333333, 0.008504666666666666, 0.025514, 0.3779124020559003, 0.3779124020559003, 0.002541545926222222, 0.022873913336, 0.1372944987953624, 0.1372944987953624),
    @("ECs", "Efna4", "Epha5", "MuSCs", 2, 0.6666666666666666, 0.2988413333333333, 0.896524, 0.3632971504731247, 0.3632971504731246, 1, 0.3333333333333333, 0.01399966666666667, 0.041999, 0.6220875979440996, 0.6220875979440996, 0.004183679052888889, 0.037653111476, 0.2260026516777622, 0.2260026516777622),
    @("FAPs", "Efna4", "Epha5", "ECs", 2, 0.6666666666666666, 0.3851916666666667, 1.155575, 0.4682720202225272, 0.4682720202225272, 1, 0.3333333333333333, 0.008504666666666666, 0.025514, 0.3779124020559003, 0.3779124020559003, 0.003275926727777777, 0.02948334055, 0.1769658039778644, 0.1769658039778644),
    @("FAPs", "Efna4", "Epha5", "MuSCs", 2, 0.6666666666666666, 0.3851916666666667, 1.155575, 0.4682720202225272, 0.4682720202225272, 1, 0.3333333333333333, 0.01399966666666667, 0.041999, 0.6220875979440996, 0.6220875979440996, 0.005392554936111111, 0.048532994425, 0.2913062162446628, 0.2913062162446628),
    @("MuSCs", "Efna4", "Epha5", "ECs", 2, 0.6666666666666666, 0.138548, 0.415644, 0.1684308293043481, 0.1684308293043481, 1, 0.3333333333333333, 0.008504666666666666, 0.025514, 0.3779124020559003, 0.3779124020559003, 0.001178304557333333, 0.010604741016, 0.06365209928267353, 0.06365209928267353),
    @("MuSCs", "Efna4", "Epha5", "MuSCs", 2, 0.6666666666666666, 0.138548, 0.415644, 0.1684308293043481, 0.1684308293043481, 1, 0.3333333333333333, 0.01399966666666667, 0.041999, 0.6220875979440996, 0.6220875979440996, 0.001939625817333334, 0.017456632356, 0.1047787300216746, 0.1047787300216746),
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $colNum = $c + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$c]
    }
}

# The sheet now spans 7 rows (header + 6 data rows) instead of the original 4.
Write-Output ("Updated rows 2-7 with new TPM data; used range now " + $ws.UsedRange.Rows.Count.ToString() + " rows x " + $ws.UsedRange.Columns.Count.ToString() + " cols")
